$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34
$ws.Range("A34").Value = 112038603
$ws.Range("B34").Value = 89503
$ws.Range("D34").Value = "LC"
$ws.Range("E34").Value = 5447
$ws.Range("F34").Value = "Vedticka"
$ws.Range("G34").Value = "Fuscoporia viticola"
$ws.Range("H34").Value = "(Schwein.) Murrill"
$ws.Range("Q34").Value = 615968
$ws.Range("R34").Value = 6895406

# Row 35
$ws.Range("A35").Value = 112038600
$ws.Range("B35").Value = 86357
$ws.Range("D35").Value = "NT"
$ws.Range("E35").Value = 4412
$ws.Range("F35").Value = "Äggvaxskivling"
$ws.Range("G35").Value = "Hygrophorus karstenii"
$ws.Range("H35").Value = "Sacc. & Cub."
$ws.Range("Q35").Value = 616034
$ws.Range("R35").Value = 6895585

# Row 36
$ws.Range("A36").Value = 112038604
$ws.Range("B36").Value = 89979
$ws.Range("D36").Value = "VU"
$ws.Range("E36").Value = 1209
$ws.Range("F36").Value = "Rynkskinn"
$ws.Range("G36").Value = "Phlebia centrifuga"
$ws.Range("H36").Value = "P.Karst."
$ws.Range("Q36").Value = 615978
$ws.Range("R36").Value = 6895550

# Row 37
$ws.Range("A37").Value = 112038601
$ws.Range("B37").Value = 73758
$ws.Range("D37").Value = "LC"
$ws.Range("E37").Value = 6426
$ws.Range("F37").Value = "Kattfotslav"
$ws.Range("G37").Value = "Felipes leucopellaeus"
$ws.Range("H37").Value = "(Ach.) Frisch & G.Thor"
$ws.Range("Q37").Value = 616013
$ws.Range("R37").Value = 6895612

# Row 38
$ws.Range("A38").Value = 112038602
$ws.Range("B38").Value = 86357
$ws.Range("D38").Value = "NT"
$ws.Range("E38").Value = 4412
$ws.Range("F38").Value = "Äggvaxskivling"
$ws.Range("G38").Value = "Hygrophorus karstenii"
$ws.Range("H38").Value = "Sacc. & Cub."
$ws.Range("Q38").Value = 616026
$ws.Range("R38").Value = 6895554

# Row 39
$ws.Range("A39").Value = 112038599
$ws.Range("B39").Value = 89557
$ws.Range("D39").Value = "NT"
$ws.Range("E39").Value = 5432
$ws.Range("F39").Value = "Granticka"
$ws.Range("G39").Value = "Porodaedalea chrysoloma"
$ws.Range("H39").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q39").Value = 616070
$ws.Range("R39").Value = 6895500

# Row 40
$ws.Range("A40").Value = 112038596
$ws.Range("B40").Value = 90221
$ws.Range("D40").Value = "LC"
$ws.Range("E40").Value = 3298
$ws.Range("F40").Value = "Trådticka"
$ws.Range("G40").Value = "Climacocystis borealis"
$ws.Range("H40").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q40").Value = 616076
$ws.Range("R40").Value = 6895428
